# Insert a new data row (row 3) into the weekly price sheet.
# This shifts all existing rows 3-26 down to 4-27 and populates
# the new row 3 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (pushes rows 3..26 down to 4..27)
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with values
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44756
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 100112044
$ws.Range("G3").Value = "Perejil"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 104
$ws.Range("K3").Value = 2800
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2904
$ws.Range("N3").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O3").Value = "Provincia de Quillota"
$ws.Range("P3").Value = 968
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Hortaliza"
